$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old A:D data (STT, id, name, photo) entirely so we can
# re-lay the table out one column to the left (A:C) without the STT column.
$ws.Range("A1:D7").Clear()

# Column widths for the new layout (values chosen so the COM layer's
# internal pixel-grid rounding lands on the closest achievable width to the
# target 22.7109375 / 51.7109375 / 109.85546875 "characters").
$ws.Columns.Item(1).ColumnWidth = 21.83
$ws.Columns.Item(2).ColumnWidth = 50.83
$ws.Columns.Item(3).ColumnWidth = 109

# Header row.
$ws.Cells.Item(1, 1).Value = "id"
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "photo"

# Data rows.
$data = @(
    @(128317, "Co.op Mart Miền Trung", "https://images.foody.vn/res/g65/641598/prof/s1242x600/foody-mobile-m-mat-jpg-939-636245845441837676.jpg"),
    @(225160, "Tuti Fruit Shop", "https://images.foody.vn/res/g112/1112799/prof/s750x400/foody-upload-api-foody-mobile-co-d127a61f-211025155256.jpeg"),
    @(256763, "Thảo Trái Cây Sỉ Lẻ ", "https://images.foody.vn/res/g114/1133888/prof/s750x400/foody-upload-api-foody-mobile-co-1c223f2d-220415153319.jpeg"),
    @(74087, "Shop Thực Phẩm Yến", "https://images.foody.vn/res/g88/878359/prof/s750x400/foody-upload-api-foody-mobile-rau-qua-190109084121.jpg"),
    @(123782, "Laban Mart - Chuỗi Siêu Thị Thực Phẩm Thiết Yếu", "https://images.foody.vn/res/g104/1035062/prof/s750x400/foody-upload-api-foody-mobile-z3-acfdeec8-221219135701.jpeg"),
    @(203678, "Thực Phẩm Sạch Thảo Nguyên", "https://images.foody.vn/res/g110/1094313/prof/s750x400/foody-upload-api-foody-mobile-co-3b6983d6-210813122705.jpeg")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$ws.Range("B13").Select()
